$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 data for the "Corrupted Bishop" raid.
$ws.Range("A6").Value = "Corrupted Bishop"
$ws.Range("B6").Value = "I will convert you godless heathens to the light. I will put down the wicked and bathe in the blood of your children’s screams. You are nothing more then a blight on the soil of the holy lord. My knights march, my priests pray, my clerics heal the wicked. You are a corruption of all that is holy. I shall stomp you in obliteration. You are nothing child. Nothing."
$ws.Range("C6").Value = "corrupted-bishop"
$ws.Range("D6").Value = "Corrupted Bishop"
$ws.Range("E6").Value = "Delusional Soul Crusher,Demonic Infestation of The Child,Corrupted Priest of The Federation,Shadow Jester of Rage,Dancing Queen of Yesterday,Faithful Cleric o fThe Church,Bloody Knight of Horror"
$ws.Range("F6").Value = "Federation Controlled Town"
$ws.Range("G6").Value = "Federation Controlled Town,Northren Port,Southren Port,Federation City"
$ws.Range("H6").Value = "Delusional Silver"
$ws.Range("I6").Value = "Ancestral Soldiers Statue"

# The new row's text in columns C (raid_type) and E (raid_monster_ids) is
# wider than the previous best-fit column widths, so Excel recalculates
# (widens) those two columns to fit the new content.
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(5).ColumnWidth = 229.16666666666666
